# Apply the "Themed Style 1 - Accent 1" table style to the B1 "Types of
# financial documents" table on slide 5, replacing the previously-applied
# custom table style.
#
# PowerPoint table styles cannot be reassigned via the `Table.Style`
# property (it is read-only); the real COM object model exposes
# `Table.ApplyStyle(styleId)` for this purpose, so that is what we call
# here, passing the GUID of the built-in style shown in the Table Styles
# gallery.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table

$table.ApplyStyle("{FAC8DE47-0582-46F3-ABE3-895D6D16116C}")
